$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings are written back as text (matching
# the original cells inlineStr/text storage) instead of being auto-coerced
# to numbers by Excel.

$ws.Range('D2').Value = '36.910.62'
$ws.Range('E2').Value = '  -1.32%  '

$ws.Range('D3').Value = '2.005.35'
$ws.Range('E3').Value = '  -2.13%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '224.95'
$ws.Range('E5').Value = '  -2.22%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.603'
$ws.Range('E6').Value = '  -1.76%  '

$ws.Range('E7').Value = '  +0.07%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '54.38'
$ws.Range('E8').Value = '  -4.82%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.372'
$ws.Range('E9').Value = '  -3.44%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0775'
$ws.Range('E10').Value = '  -3.34%  '

$ws.Range('E11').Value = '  -5.30%  '

$ws.Range('D12').Value = '2.304.40'
$ws.Range('E12').Value = '  -2.04%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '13.91'
$ws.Range('E13').Value = '  -5.55%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '19.69'
$ws.Range('E14').Value = '  -5.40%  '

$ws.Range('E15').Value = '  -2.26%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.732'
$ws.Range('E16').Value = '  -3.57%  '

$ws.Range('D17').Value = '2.052.63'
$ws.Range('E17').Value = '  +0.29%  '

$ws.Range('D18').Value = '36.850.20'
$ws.Range('E18').Value = '  -1.18%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.30'
$ws.Range('E19').Value = '  +3.37%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.10'
$ws.Range('E20').Value = '  -2.44%  '

$ws.Range('D21').Value = '0.0₃0807'
$ws.Range('E21').Value = '  -3.08%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '221.15'
$ws.Range('E22').Value = '  -2.40%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.02%  '

$ws.Range('E24').Value = '  +1.60%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.15'
$ws.Range('E25').Value = '  -6.46%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.20'
$ws.Range('E26').Value = '  -2.81%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.99'
$ws.Range('E27').Value = '  -6.22%  '

$ws.Range('E28').Value = '  -3.45%  '

$ws.Range('E29').Value = '  -2.91%  '

$ws.Range('E30').Value = '  -6.29%  '

$ws.Range('E31').Value = '  -2.38%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.43'
$ws.Range('E32').Value = '  -2.69%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0596'
$ws.Range('E33').Value = '  -3.29%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.41'
$ws.Range('E34').Value = '  -3.95%  '

$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.87'
$ws.Range('E35').Value = '  +2.02%  '

$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.30'
$ws.Range('E36').Value = '  -5.50%  '

$ws.Range('E37').Value = '  +0.06%  '

$ws.Range('E38').Value = '  -5.12%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.33'
$ws.Range('E39').Value = '  -0.77%  '

$ws.Range('D40').Value = '1.450.36'
$ws.Range('E40').Value = '  -3.07%  '

$ws.Range('E41').Value = '  -5.40%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '94.25'
$ws.Range('E42').Value = '  -2.51%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.77'
$ws.Range('E43').Value = '  -4.63%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0908'
$ws.Range('E44').Value = '  -3.50%  '

$ws.Range('E45').Value = '  -4.68%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '15.81'
$ws.Range('E46').Value = '  -8.56%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.10'
$ws.Range('E47').Value = '  -0.87%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.992'
$ws.Range('E48').Value = '  -3.01%  '

$ws.Range('E49').Value = '  -0.62%  '

$ws.Range('D50').Value = '2.192.93'
$ws.Range('E50').Value = '  -2.06%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.51'
$ws.Range('E51').Value = '  -10.40%  '
